$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new product row for "DIAMICRON 60MG M.R. 30 SCORED TAB" right
# before "DICLAC 150 ID 20 PROLONGED R TABS" (alphabetical order), which
# currently lives on row 32. This pushes every row below it down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).Insert()

# Copy the formatting (styles, borders, fonts, fill, number formats, etc.)
# of the row that used to be 32 (now shifted to row 33) into the freshly
# inserted, still-blank row 32 so the new row looks identical to the rest
# of the table.
$ws.Range("A33:Q33").Copy()
$ws.Range("A32:Q32").PasteSpecial(-4122)
$ws.Range("Q33").Copy()
$ws.Range("Q32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A32:Q32").RowHeight = 25.5

# Recreate the per-row merged cell regions (A:B, C:G, H:K, L:M, N:O) that
# every data row in the table has.
$ws.Range("A32:B32").Merge()
$ws.Range("C32:G32").Merge()
$ws.Range("H32:K32").Merge()
$ws.Range("L32:M32").Merge()
$ws.Range("N32:O32").Merge()

# Every text-like column is stored as plain text in this workbook (even the
# ones that look numeric, such as the balance count or the selling price),
# so force a "@" (text) number format before assigning the value - otherwise
# the COM layer would silently coerce a numeric-looking string into a real
# number. The final, column-specific number format is (re)applied afterwards
# purely for cosmetic/display purposes; it does not turn the cell back into
# a number once the value is already stored as text.
$ws.Range("A32").NumberFormat = "General"
$ws.Range("A32").Value = 26

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "DIAMICRON 60MG M.R. 30 SCORED TAB"

$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "3:1"

$ws.Range("L32").NumberFormat = "@"
$ws.Range("L32").Value = "1"
$ws.Range("L32").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N32").NumberFormat = "@"
$ws.Range("N32").Value = "156.00"

$ws.Range("P32").NumberFormat = "@"
$ws.Range("P32").Value = "156.0000"
$ws.Range("P32").NumberFormat = "0.00"

$ws.Range("Q32").NumberFormat = "@"
$ws.Range("Q32").Value = "0:1"

# ---------------------------------------------------------------------------
# The running-total row (previously row 103, now row 104 after the insert)
# needs to reflect the extra selling price of the new item.
# ---------------------------------------------------------------------------
$total = $ws.Range("P104").Value + 156
$ws.Range("P104").Value = $total

# ---------------------------------------------------------------------------
# Update the generation timestamp shown in the report footer (now row 105).
# ---------------------------------------------------------------------------
$ws.Range("A105").Value = "Thursday, 7 August, 2025 8:32 PM"
